# Type the document's text content into cell A1 of the active sheet.
# (The leading character is the BOM/ZWNBSP codepoint U+FEFF that was present
# at the start of the text in the source document.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text = [char]0xFEFF + "This is the content of the document."
$ws.Range("A1").Value = $text
